# Updated cryptos list on Tue Apr 11 05:39:20 UTC 2023 with GitHub Actions
#
# For each row, column D holds the "Price" and column E holds the
# "Volume(1h)" percentage text. Many Price values look like plain decimals
# (e.g. "327.05"), so a bare .Value assignment would be auto-coerced to a
# number by Excel and lose the original text formatting/trailing zeros.
# We force those through as literal text (leading apostrophe) and then
# strip the resulting "Text" cell format back off with ClearFormats() so
# the cell keeps its original (default) style, matching the source data
# which was authored as plain inline strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Numeric-looking strings (plain decimals) get auto-converted to
    # numbers by Excel unless explicitly quoted as text; values that are
    # already unambiguous text (e.g. contain two dots) don't need it but
    # quoting them is harmless since Excel strips a non-meaningful prefix.
    $range.Value = "'" + $value
    $range.ClearFormats()
}

$rows = @(
    @{ Row = 2;  D = "30.110.18";  E = "  +5.85%  " },
    @{ Row = 3;  D = "1.923.29";   E = "  +2.75%  " },
    @{ Row = 4;  D = $null;        E = "  -0.79%  " },
    @{ Row = 5;  D = "327.05";     E = "  +3.64%  " },
    @{ Row = 6;  D = $null;        E = "  -0.72%  " },
    @{ Row = 7;  D = "0.5209";     E = "  +2.72%  " },
    @{ Row = 8;  D = "0.4071";     E = "  +4.49%  " },
    @{ Row = 9;  D = "0.08484";    E = "  +1.41%  " },
    @{ Row = 10; D = "1.131";      E = "  +2.77%  " },
    @{ Row = 11; D = "42.86";      E = "  +2.81%  " },
    @{ Row = 12; D = "22.22";      E = "  +9.19%  " },
    @{ Row = 13; D = "6.373";      E = "  +2.62%  " },
    @{ Row = 14; D = "1.924.44";   E = "  +2.64%  " },
    @{ Row = 15; D = "7.386";      E = "  +2.04%  " },
    @{ Row = 16; D = $null;        E = "  -0.83%  " },
    @{ Row = 17; D = "96.14";      E = "  +5.48%  " },
    @{ Row = 18; D = $null;        E = "  +1.52%  " },
    @{ Row = 19; D = "0.06745";    E = "  +0.29%  " },
    @{ Row = 20; D = $null;        E = "  +3.41%  " },
    @{ Row = 21; D = "1.001";      E = "  -0.67%  " },
    @{ Row = 22; D = "6.071";      E = "  +2.56%  " },
    @{ Row = 23; D = "30.116.73";  E = "  +5.75%  " },
    @{ Row = 24; D = $null;        E = "  +2.04%  " },
    @{ Row = 25; D = $null;        E = "  -1.31%  " },
    @{ Row = 26; D = "2.145.55";   E = "  +2.76%  " },
    @{ Row = 27; D = "21.23";      E = "  +3.07%  " },
    @{ Row = 28; D = "160.57";     E = "  -0.64%  " },
    @{ Row = 29; D = "2.470";      E = "  +3.76%  " },
    @{ Row = 30; D = "129.03";     E = "  +2.59%  " },
    @{ Row = 31; D = "1.087";      E = "  +5.01%  " },
    @{ Row = 32; D = "0.1061";     E = "  +1.74%  " },
    @{ Row = 33; D = "6.112";      E = "  +6.14%  " },
    @{ Row = 34; D = "3.657";      E = "  +1.46%  " },
    @{ Row = 35; D = "0.02526";    E = "  +3.11%  " },
    @{ Row = 36; D = "0.06630";    E = "  +1.52%  " },
    @{ Row = 37; D = "0.2224";     E = "  +3.13%  " },
    @{ Row = 38; D = "1.242";      E = "  +4.71%  " },
    @{ Row = 39; D = "9.048";      E = "  +2.50%  " },
    @{ Row = 40; D = "5.218";      E = "  +3.18%  " },
    @{ Row = 41; D = "0.6594";     E = "  +3.25%  " },
    @{ Row = 42; D = "1.252";      E = "  -0.10%  " },
    @{ Row = 43; D = "11.64";      E = "  +4.93%  " },
    @{ Row = 44; D = "0.6199";     E = "  +3.13%  " },
    @{ Row = 45; D = "13.22";      E = "  +1.66%  " },
    @{ Row = 46; D = "3.763";      E = "  +2.26%  " },
    @{ Row = 47; D = "2.078";      E = "  +3.81%  " },
    @{ Row = 48; D = "1.247";      E = "  +2.89%  " },
    @{ Row = 49; D = "125.95";     E = "  +3.57%  " },
    @{ Row = 50; D = "1.165";      E = "  +3.33%  " },
    @{ Row = 51; D = $null;        E = "  +4.55%  " }
)

foreach ($r in $rows) {
    if ($null -ne $r.D) {
        Set-TextValue $ws.Range("D" + $r.Row) $r.D
    }
    $ws.Range("E" + $r.Row).Value = $r.E
}
